# "Save HSI, option for core temp/HSI depending on HSI range, some optimization for memory"
#
# The "Configure Fellrnr DataField" sheet lets the user pick abbreviated
# field codes (via a Lookup-sheet driven dropdown) that get concatenated
# into copy/paste strings (A2/A3/A4, built from K7/K16/K24 respectively).
# Update a few of those dropdown picks:
#   D8  / D17 : cadence ("ca")         -> dist ("dis")
#   D25       : cadence ("ca")         -> deltaElevation ("dEl")
#   D11       : averagePace ("aP")     -> heat ("ht")  (HSI = Heat Stress Index)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Configure Fellrnr DataField")

$ws1.Range("D8").Value = "dist"
$ws1.Range("D17").Value = "dist"
$ws1.Range("D25").Value = "deltaElevation"
$ws1.Range("D11").Value = "heat"

# Bring the "Configure Fellrnr DataField" sheet to the front (it was the
# "Lookup" sheet that was active/selected before) with D11 as the active
# selected cell, matching the saved workbook view state.
$ws1.Activate()
$ws1.Range("D11").Select()
